# ooutput update 2025 august
#
# Refreshes the generated IG "output" workbook for the new publication run:
#   - the canonical URL base moved from the old GitHub shorthand mirror to
#     the new 2rdoc.pt IG site (StructureDefinition + ValueSet links)
#   - the IG "Date" metadata value was regenerated for the August run
#   - the "Elements" sheet's best-fit column widths shift slightly because
#     the regenerated cell content (URLs/date) is a different length than
#     before; widths are refreshed to match the new best-fit layout

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Metadata sheet: canonical URL + publication Date
# ---------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/StructureDefinition/stress-triggers"
$wsMeta.Range("B8").Value = "2025-08-20T10:40:04+01:00"

# ---------------------------------------------------------------------
# Elements sheet: Binding Value Set URL
# ---------------------------------------------------------------------
$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("Z6").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/ValueSet/stress-triggers-vs"

# ---------------------------------------------------------------------
# Elements sheet: refreshed best-fit column widths (regenerated layout).
# ColumnWidth is specified in character units; Excel persists it back to
# the sheet's stored column width (character width + the fixed ~5px cell
# padding), so the values below are chosen so the stored/persisted width
# matches the regenerated best-fit widths.
# ---------------------------------------------------------------------
$wsElem.Columns.Item(1).ColumnWidth = 15.666666666666666
$wsElem.Columns.Item(2).ColumnWidth = 15.666666666666666
$wsElem.Columns.Item(3).ColumnWidth = 9.0
$wsElem.Columns.Item(4).ColumnWidth = 6.166666666666667
$wsElem.Columns.Item(5).ColumnWidth = 4.5
$wsElem.Columns.Item(6).ColumnWidth = 3.1666666666666665
$wsElem.Columns.Item(7).ColumnWidth = 3.5
$wsElem.Columns.Item(8).ColumnWidth = 11.833333333333334
$wsElem.Columns.Item(9).ColumnWidth = 9.666666666666666
$wsElem.Columns.Item(11).ColumnWidth = 13.5
$wsElem.Columns.Item(15).ColumnWidth = 11.5
$wsElem.Columns.Item(20).ColumnWidth = 7.0
$wsElem.Columns.Item(21).ColumnWidth = 12.833333333333334
$wsElem.Columns.Item(22).ColumnWidth = 13.166666666666666
$wsElem.Columns.Item(23).ColumnWidth = 14.166666666666666
$wsElem.Columns.Item(24).ColumnWidth = 13.833333333333334
$wsElem.Columns.Item(25).ColumnWidth = 16.166666666666668
$wsElem.Columns.Item(26).ColumnWidth = 55.0
$wsElem.Columns.Item(27).ColumnWidth = 4.166666666666667
$wsElem.Columns.Item(28).ColumnWidth = 17.166666666666668
$wsElem.Columns.Item(29).ColumnWidth = 33.666666666666664
$wsElem.Columns.Item(30).ColumnWidth = 12.666666666666666
$wsElem.Columns.Item(31).ColumnWidth = 10.5
$wsElem.Columns.Item(32).ColumnWidth = 14.166666666666666
$wsElem.Columns.Item(33).ColumnWidth = 7.333333333333333
$wsElem.Columns.Item(34).ColumnWidth = 7.666666666666667
$wsElem.Columns.Item(37).ColumnWidth = 18.666666666666668
